$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Turn the empty paragraph right after the title "Position Paper" into a
#    centered, bold "Author: Dr Kevin R Bond" line.
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleFound = $titleRange.Find.Execute("Position Paper", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($titleFound) {
    $titlePara   = $titleRange.Paragraphs(1)
    $authorPara  = $titlePara.Next()
    $authorRange = $authorPara.Range

    $authorRange.Text = "Author: Dr Kevin R Bond"
    $authorRange.Font.Bold = 1
    $authorPara.Alignment = 1   # wdAlignParagraphCenter
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from its old spot (near the end of the
#    document) to the empty paragraph right after the "... assembly
#    languages." list item (the paragraph whose only pPr content is
#    <w:ind w:left="1080"/>).
# ---------------------------------------------------------------------------
$goBack = $null
try {
    $goBack = $d.Bookmarks("_GoBack")
} catch {
    $goBack = $null
}
if ($goBack -ne $null) {
    $goBack.Delete()
}

$langRange = $d.Content
$langFound = $langRange.Find.Execute("Lisp, VHDL, Verilog, assembly languages", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($langFound) {
    $langPara    = $langRange.Paragraphs(1)
    $targetPara  = $langPara.Next()
    $targetRange = $targetPara.Range
    $d.Bookmarks.Add("_GoBack", $targetRange)
}
